# This script applies the following edit to the document:
#   1. Remove the paragraph containing "123456" entirely.
#   2. Change the "2018020054" paragraph's paragraph-mark rFonts hint
#      from "default" (w:eastAsiaTheme="minorEastAsia") to "eastAsia".
#   3. Insert a brand-new paragraph ("2018020118") right after it, moving
#      the "_GoBack" bookmark from the "2018020054" paragraph onto this
#      new, final paragraph.

$d = $word.ActiveDocument

# --- Step 1: delete the paragraph whose text is "123456" (whole paragraph,
#     including its end-of-paragraph mark) ---
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "123456`r") {
        $p.Range.Delete()
    }
}

# --- Step 2 & 3: locate the "2018020054" paragraph and replace it (text +
#     end-of-paragraph mark) with the updated "2018020054" paragraph
#     (rFonts hint now "eastAsia") immediately followed by the new
#     "2018020118" paragraph carrying the "_GoBack" bookmark. ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "2018020054`r") {
        $targetRange = $d.Range($p.Range.Start, $p.Range.End)

        $wordml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
      <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
      <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
    </w:rPr>
    <w:t>2018020054</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:hint="default"/>
      <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:hint="eastAsia"/>
      <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
    </w:rPr>
    <w:t>2018020118</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@

        $targetRange.InsertXML($wordml)
        break
    }
}
